$d = $word.ActiveDocument

function Get-ParagraphContaining($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# 1. Mark these LeetCode problems as "finished": red text (FF0000) + single
#    underline, matching the styling already used for #1, #17, #20, #79, #206...
$finishedNeedles = @(
    "Remove Duplicates from Sorted Array II",
    "Validate Binary Search Tree",
    "Word Ladder",
    "Clone Graph"
)
foreach ($needle in $finishedNeedles) {
    $p = Get-ParagraphContaining $d $needle
    $r = $p.Range
    $r.Font.Color = 255
    $r.Font.Underline = 1
}

# 2. Paragraph 139 "Word Break": the cursor (_GoBack) now lives around the
#    "Word Break" text, splitting the single run into three runs.
$p139 = Get-ParagraphContaining $d "Word Break"
$text139 = $p139.Range.Text
$wbStart = $text139.IndexOf("Word Break")
$wbEnd = $wbStart + "Word Break        ".Length
$goBackRange = $d.Range($p139.Range.Start + $wbStart, $p139.Range.Start + $wbEnd)
# Bookmark names are unique, so adding "_GoBack" here removes it from wherever
# it used to be (paragraph 334, "Increasing Triplet Subsequence").
$d.Bookmarks.Add("_GoBack", $goBackRange)

# 3. Paragraph 334 "...Increasing Triplet Subsequence...": losing its _GoBack
#    bookmark leaves the text split across two runs ("334        In" /
#    "creasing..."); re-merge them into a single run with a no-op replace that
#    spans the old run boundary.
$p334 = Get-ParagraphContaining $d "Increasing Triplet Subsequence"
$p334.Range.Find.Execute("Increasing", $true, $false, $false, $false, $false, $true, 1, $false, "Increasing", 2)

Write-Output "done"
